{"js": "// Change \"REFER TO THE FOLLOWING FOR AN EXHAUSTIVE LIST:\" to\n// \"REFER TO THE FOLLOWING FOR A MORE EXHAUSTIVE LIST:\" and move the\n// auto-managed \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n// to the new edit location.\n\nconst body = context.document.body;\n\n// Word only ever keeps a single \"_GoBack\" bookmark \u2014 editing text moves\n// it to the new edit location. Remove the old one first (wherever it is).\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no existing _GoBack bookmark - nothing to remove\n}\n\n// Step 1: locate the \"N\" in \"...FOR AN EXHAUSTIVE LIST:\" and delete it,\n// turning \"AN\" into \"A\".\nconst tailSearch = body.search(\"N EXHAUSTIVE LIST:\", { matchCase: true });\ntailSearch.load(\"text\");\nawait context.sync();\nconst tail = tailSearch.items[0];\n\nconst nSearch = tail.search(\"N\", { matchCase: true });\nnSearch.load(\"text\");\nawait context.sync();\nconst nChar = nSearch.items[0];\n\nnChar.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: insert \" MORE\" right before \" EXHAUSTIVE LIST:\".\nconst exhSearch = body.search(\" EXHAUSTIVE LIST:\", { matchCase: true });\nexhSearch.load(\"text\");\nawait context.sync();\nconst exh = exhSearch.items[0];\n\nconst insertedRange = exh.insertText(\" MORE\", Word.InsertLocation.before);\nawait context.sync();\n\n// Step 3: drop the \"_GoBack\" bookmark at the new edit location (right\n// after the just-inserted \" MORE\" text, before \" EXHAUSTIVE LIST:\").\ninsertedRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst moreSearch = body.search(\" MORE\", { matchCase: true });\nmoreSearch.load(\"text\");\nawait context.sync();\nconst moreRange = moreSearch.items[0];\nconst afterMore = moreRange.getRange(\"End\");\nafterMore.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change \"REFER TO THE FOLLOWING FOR AN EXHAUSTIVE LIST:\" to\n# \"REFER TO THE FOLLOWING FOR A MORE EXHAUSTIVE LIST:\" and move the\n# auto-managed \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n# to the new edit location.\n\n$d = $word.ActiveDocument\n\n# Step 1: locate the \"N\" in \"...FOR AN EXHAUSTIVE LIST:\" and delete it,\n# turning \"AN\" into \"A\".\n$tailRange = $d.Content\n$tailRange.Find.Text = \"N EXHAUSTIVE LIST:\"\n$tailRange.Find.Execute()\n$nStart = $tailRange.Start\n$nRange = $d.Range($nStart, $nStart + 1)\n$nRange.Text = \"\"\n\n# Step 2: insert \" MORE\" right before \" EXHAUSTIVE LIST:\".\n$exhRange = $d.Content\n$exhRange.Find.Text = \" EXHAUSTIVE LIST:\"\n$exhRange.Find.Execute()\n$exhStart = $exhRange.Start\n$insertPoint = $d.Range($exhStart, $exhStart)\n$insertPoint.InsertBefore(\" MORE\")\n\n# Step 3: drop the \"_GoBack\" bookmark at the new edit location (right\n# after the just-inserted \" MORE\" text, before \" EXHAUSTIVE LIST:\").\n# Word keeps only a single \"_GoBack\" bookmark, so adding it here moves\n# it from wherever it previously was (the end of the document).\n#\n# Bookmarking the \" MORE\" span first (instead of going straight to a\n# collapsed bookmark) keeps \"...FOR A\" and \" MORE\" as distinct runs,\n# matching how Word itself preserves run boundaries around edits.\n$moreRange = $d.Content\n$moreRange.Find.Text = \" MORE\"\n$moreRange.Find.Execute()\n$d.Bookmarks.Add(\"_GoBack\", $moreRange)\n\n$moreRange2 = $d.Content\n$moreRange2.Find.Text = \" MORE\"\n$moreRange2.Find.Execute()\n$afterMoreStart = $moreRange2.End\n$collapsed = $d.Range($afterMoreStart, $afterMoreStart)\n$d.Bookmarks.Add(\"_GoBack\", $collapsed)\n"}
